{"js": "// Positional update of the date paragraph and the 100 multiplication-table\n// cells. The document body contains exactly 101 paragraphs in this fixed\n// order: the date line, then the 100 \"<a>x<b>=\" table-cell paragraphs\n// (row-major, top-to-bottom / left-to-right). Each entry below is\n// [oldText, newText] for one paragraph, in that same order.\nconst replacements = [\n  [\"2023-03-21 Tuesday\", \"2023-03-22 Wednesday\"],\n  [\"50\u00d760=\", \"23\u00d761=\"],\n  [\"20\u00d763=\", \"55\u00d755=\"],\n  [\"77\u00d798=\", \"72\u00d786=\"],\n  [\"37\u00d710=\", \"19\u00d795=\"],\n  [\"30\u00d799=\", \"96\u00d770=\"],\n  [\"82\u00d777=\", \"64\u00d784=\"],\n  [\"63\u00d722=\", \"65\u00d757=\"],\n  [\"28\u00d799=\", \"86\u00d794=\"],\n  [\"32\u00d719=\", \"10\u00d766=\"],\n  [\"45\u00d766=\", \"16\u00d759=\"],\n  [\"38\u00d790=\", \"17\u00d748=\"],\n  [\"69\u00d787=\", \"67\u00d760=\"],\n  [\"19\u00d799=\", \"14\u00d7100=\"],\n  [\"95\u00d793=\", \"36\u00d722=\"],\n  [\"98\u00d752=\", \"34\u00d754=\"],\n  [\"43\u00d712=\", \"44\u00d710=\"],\n  [\"92\u00d748=\", \"99\u00d760=\"],\n  [\"55\u00d725=\", \"60\u00d735=\"],\n  [\"85\u00d745=\", \"100\u00d787=\"],\n  [\"47\u00d769=\", \"37\u00d720=\"],\n  [\"99\u00d759=\", \"32\u00d711=\"],\n  [\"50\u00d736=\", \"70\u00d737=\"],\n  [\"28\u00d771=\", \"25\u00d724=\"],\n  [\"94\u00d718=\", \"90\u00d730=\"],\n  [\"54\u00d790=\", \"45\u00d734=\"],\n  [\"96\u00d766=\", \"80\u00d758=\"],\n  [\"46\u00d799=\", \"17\u00d735=\"],\n  [\"24\u00d793=\", \"39\u00d716=\"],\n  [\"47\u00d741=\", \"42\u00d716=\"],\n  [\"23\u00d711=\", \"13\u00d736=\"],\n  [\"53\u00d766=\", \"23\u00d747=\"],\n  [\"49\u00d768=\", \"45\u00d762=\"],\n  [\"45\u00d725=\", \"77\u00d742=\"],\n  [\"84\u00d773=\", \"63\u00d742=\"],\n  [\"77\u00d741=\", \"67\u00d775=\"],\n  [\"94\u00d784=\", \"54\u00d741=\"],\n  [\"25\u00d729=\", \"75\u00d741=\"],\n  [\"72\u00d727=\", \"15\u00d763=\"],\n  [\"46\u00d769=\", \"54\u00d761=\"],\n  [\"99\u00d712=\", \"77\u00d789=\"],\n  [\"12\u00d791=\", \"41\u00d715=\"],\n  [\"22\u00d774=\", \"53\u00d716=\"],\n  [\"66\u00d787=\", \"20\u00d727=\"],\n  [\"97\u00d743=\", \"13\u00d746=\"],\n  [\"13\u00d730=\", \"20\u00d726=\"],\n  [\"82\u00d779=\", \"57\u00d763=\"],\n  [\"85\u00d787=\", \"29\u00d775=\"],\n  [\"83\u00d789=\", \"99\u00d743=\"],\n  [\"81\u00d774=\", \"47\u00d710=\"],\n  [\"47\u00d722=\", \"95\u00d727=\"],\n  [\"96\u00d762=\", \"26\u00d735=\"],\n  [\"92\u00d727=\", \"38\u00d733=\"],\n  [\"20\u00d778=\", \"46\u00d766=\"],\n  [\"31\u00d766=\", \"30\u00d724=\"],\n  [\"61\u00d736=\", \"14\u00d735=\"],\n  [\"58\u00d772=\", \"41\u00d723=\"],\n  [\"74\u00d789=\", \"22\u00d760=\"],\n  [\"89\u00d736=\", \"10\u00d7100=\"],\n  [\"92\u00d761=\", \"78\u00d795=\"],\n  [\"30\u00d733=\", \"65\u00d778=\"],\n  [\"67\u00d752=\", \"90\u00d780=\"],\n  [\"59\u00d764=\", \"83\u00d782=\"],\n  [\"11\u00d724=\", \"55\u00d718=\"],\n  [\"74\u00d749=\", \"79\u00d751=\"],\n  [\"93\u00d765=\", \"65\u00d744=\"],\n  [\"59\u00d749=\", \"20\u00d773=\"],\n  [\"41\u00d742=\", \"82\u00d745=\"],\n  [\"65\u00d741=\", \"39\u00d798=\"],\n  [\"59\u00d781=\", \"76\u00d769=\"],\n  [\"96\u00d734=\", \"42\u00d760=\"],\n  [\"93\u00d728=\", \"60\u00d770=\"],\n  [\"63\u00d797=\", \"43\u00d727=\"],\n  [\"85\u00d712=\", \"42\u00d789=\"],\n  [\"16\u00d743=\", \"91\u00d786=\"],\n  [\"96\u00d794=\", \"74\u00d771=\"],\n  [\"13\u00d738=\", \"84\u00d737=\"],\n  [\"80\u00d716=\", \"65\u00d755=\"],\n  [\"14\u00d760=\", \"32\u00d783=\"],\n  [\"26\u00d743=\", \"38\u00d712=\"],\n  [\"70\u00d713=\", \"50\u00d783=\"],\n  [\"48\u00d728=\", \"62\u00d764=\"],\n  [\"72\u00d715=\", \"55\u00d773=\"],\n  [\"60\u00d722=\", \"97\u00d777=\"],\n  [\"28\u00d745=\", \"77\u00d730=\"],\n  [\"73\u00d7100=\", \"36\u00d745=\"],\n  [\"82\u00d741=\", \"16\u00d722=\"],\n  [\"11\u00d718=\", \"13\u00d710=\"],\n  [\"19\u00d777=\", \"33\u00d737=\"],\n  [\"22\u00d781=\", \"43\u00d797=\"],\n  [\"56\u00d728=\", \"58\u00d782=\"],\n  [\"34\u00d748=\", \"96\u00d763=\"],\n  [\"88\u00d798=\", \"87\u00d716=\"],\n  [\"43\u00d780=\", \"41\u00d736=\"],\n  [\"38\u00d763=\", \"31\u00d792=\"],\n  [\"92\u00d734=\", \"48\u00d737=\"],\n  [\"57\u00d787=\", \"71\u00d755=\"],\n  [\"35\u00d764=\", \"62\u00d718=\"],\n  [\"29\u00d717=\", \"97\u00d785=\"],\n  [\"58\u00d743=\", \"61\u00d798=\"],\n  [\"82\u00d775=\", \"38\u00d794=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\n    \"Expected \" + replacements.length + \" paragraphs, found \" +\n    paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = paragraphs.items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \": expected \" + JSON.stringify(oldText) +\n      \" but found \" + JSON.stringify(para.text)\n    );\n  }\n  para.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Positional/unique text replacement for the date line and the 100\n# \"<a>x<b>=\" multiplication-table cells. Every old value below is unique\n# in the document, so Find/Execute from the top of the content each time\n# locates exactly the intended run without disturbing its formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-03-21 Tuesday\", \"2023-03-22 Wednesday\"),\n    @(\"50\u00d760=\", \"23\u00d761=\"),\n    @(\"20\u00d763=\", \"55\u00d755=\"),\n    @(\"77\u00d798=\", \"72\u00d786=\"),\n    @(\"37\u00d710=\", \"19\u00d795=\"),\n    @(\"30\u00d799=\", \"96\u00d770=\"),\n    @(\"82\u00d777=\", \"64\u00d784=\"),\n    @(\"63\u00d722=\", \"65\u00d757=\"),\n    @(\"28\u00d799=\", \"86\u00d794=\"),\n    @(\"32\u00d719=\", \"10\u00d766=\"),\n    @(\"45\u00d766=\", \"16\u00d759=\"),\n    @(\"38\u00d790=\", \"17\u00d748=\"),\n    @(\"69\u00d787=\", \"67\u00d760=\"),\n    @(\"19\u00d799=\", \"14\u00d7100=\"),\n    @(\"95\u00d793=\", \"36\u00d722=\"),\n    @(\"98\u00d752=\", \"34\u00d754=\"),\n    @(\"43\u00d712=\", \"44\u00d710=\"),\n    @(\"92\u00d748=\", \"99\u00d760=\"),\n    @(\"55\u00d725=\", \"60\u00d735=\"),\n    @(\"85\u00d745=\", \"100\u00d787=\"),\n    @(\"47\u00d769=\", \"37\u00d720=\"),\n    @(\"99\u00d759=\", \"32\u00d711=\"),\n    @(\"50\u00d736=\", \"70\u00d737=\"),\n    @(\"28\u00d771=\", \"25\u00d724=\"),\n    @(\"94\u00d718=\", \"90\u00d730=\"),\n    @(\"54\u00d790=\", \"45\u00d734=\"),\n    @(\"96\u00d766=\", \"80\u00d758=\"),\n    @(\"46\u00d799=\", \"17\u00d735=\"),\n    @(\"24\u00d793=\", \"39\u00d716=\"),\n    @(\"47\u00d741=\", \"42\u00d716=\"),\n    @(\"23\u00d711=\", \"13\u00d736=\"),\n    @(\"53\u00d766=\", \"23\u00d747=\"),\n    @(\"49\u00d768=\", \"45\u00d762=\"),\n    @(\"45\u00d725=\", \"77\u00d742=\"),\n    @(\"84\u00d773=\", \"63\u00d742=\"),\n    @(\"77\u00d741=\", \"67\u00d775=\"),\n    @(\"94\u00d784=\", \"54\u00d741=\"),\n    @(\"25\u00d729=\", \"75\u00d741=\"),\n    @(\"72\u00d727=\", \"15\u00d763=\"),\n    @(\"46\u00d769=\", \"54\u00d761=\"),\n    @(\"99\u00d712=\", \"77\u00d789=\"),\n    @(\"12\u00d791=\", \"41\u00d715=\"),\n    @(\"22\u00d774=\", \"53\u00d716=\"),\n    @(\"66\u00d787=\", \"20\u00d727=\"),\n    @(\"97\u00d743=\", \"13\u00d746=\"),\n    @(\"13\u00d730=\", \"20\u00d726=\"),\n    @(\"82\u00d779=\", \"57\u00d763=\"),\n    @(\"85\u00d787=\", \"29\u00d775=\"),\n    @(\"83\u00d789=\", \"99\u00d743=\"),\n    @(\"81\u00d774=\", \"47\u00d710=\"),\n    @(\"47\u00d722=\", \"95\u00d727=\"),\n    @(\"96\u00d762=\", \"26\u00d735=\"),\n    @(\"92\u00d727=\", \"38\u00d733=\"),\n    @(\"20\u00d778=\", \"46\u00d766=\"),\n    @(\"31\u00d766=\", \"30\u00d724=\"),\n    @(\"61\u00d736=\", \"14\u00d735=\"),\n    @(\"58\u00d772=\", \"41\u00d723=\"),\n    @(\"74\u00d789=\", \"22\u00d760=\"),\n    @(\"89\u00d736=\", \"10\u00d7100=\"),\n    @(\"92\u00d761=\", \"78\u00d795=\"),\n    @(\"30\u00d733=\", \"65\u00d778=\"),\n    @(\"67\u00d752=\", \"90\u00d780=\"),\n    @(\"59\u00d764=\", \"83\u00d782=\"),\n    @(\"11\u00d724=\", \"55\u00d718=\"),\n    @(\"74\u00d749=\", \"79\u00d751=\"),\n    @(\"93\u00d765=\", \"65\u00d744=\"),\n    @(\"59\u00d749=\", \"20\u00d773=\"),\n    @(\"41\u00d742=\", \"82\u00d745=\"),\n    @(\"65\u00d741=\", \"39\u00d798=\"),\n    @(\"59\u00d781=\", \"76\u00d769=\"),\n    @(\"96\u00d734=\", \"42\u00d760=\"),\n    @(\"93\u00d728=\", \"60\u00d770=\"),\n    @(\"63\u00d797=\", \"43\u00d727=\"),\n    @(\"85\u00d712=\", \"42\u00d789=\"),\n    @(\"16\u00d743=\", \"91\u00d786=\"),\n    @(\"96\u00d794=\", \"74\u00d771=\"),\n    @(\"13\u00d738=\", \"84\u00d737=\"),\n    @(\"80\u00d716=\", \"65\u00d755=\"),\n    @(\"14\u00d760=\", \"32\u00d783=\"),\n    @(\"26\u00d743=\", \"38\u00d712=\"),\n    @(\"70\u00d713=\", \"50\u00d783=\"),\n    @(\"48\u00d728=\", \"62\u00d764=\"),\n    @(\"72\u00d715=\", \"55\u00d773=\"),\n    @(\"60\u00d722=\", \"97\u00d777=\"),\n    @(\"28\u00d745=\", \"77\u00d730=\"),\n    @(\"73\u00d7100=\", \"36\u00d745=\"),\n    @(\"82\u00d741=\", \"16\u00d722=\"),\n    @(\"11\u00d718=\", \"13\u00d710=\"),\n    @(\"19\u00d777=\", \"33\u00d737=\"),\n    @(\"22\u00d781=\", \"43\u00d797=\"),\n    @(\"56\u00d728=\", \"58\u00d782=\"),\n    @(\"34\u00d748=\", \"96\u00d763=\"),\n    @(\"88\u00d798=\", \"87\u00d716=\"),\n    @(\"43\u00d780=\", \"41\u00d736=\"),\n    @(\"38\u00d763=\", \"31\u00d792=\"),\n    @(\"92\u00d734=\", \"48\u00d737=\"),\n    @(\"57\u00d787=\", \"71\u00d755=\"),\n    @(\"35\u00d764=\", \"62\u00d718=\"),\n    @(\"29\u00d717=\", \"97\u00d785=\"),\n    @(\"58\u00d743=\", \"61\u00d798=\"),\n    @(\"82\u00d775=\", \"38\u00d794=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $oldText, $false, $false, $false, $false, $false,\n        $true, 1, $false, $newText, 2\n    )\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n\n"}
